$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column B to fit the new, longer hint strings ---
$ws.Columns("B").ColumnWidth = 43.44140625

# --- Row 8: "Product Page check for hints regarding mandatory fields" ---
# Column A of this "test description" row keeps the default/general style,
# just like rows 6 and 7, so copy formats from A6 (no explicit style) first.
$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Product Page check for hints regarding mandatory fields"

$ws.Range("C8").Value = "<HINT Select an option>"
$ws.Range("D8").Value = "<HINT Select an option>"
$ws.Range("E8").Value = "<HINT Select an option>"
$ws.Range("I8").Value = "<HINT Select an option>"
$ws.Range("B8").Value = "<HINT This field is mandatory>"
$ws.Range("F8").Value = "<HINT Select at least 1 options>"

# --- Row 9: "Product Page enter date with invalid format" ---
$ws.Range("B9").Value = "0815"

# --- Row 10: "Product Page check for hint date with invalid format" ---
$ws.Range("A6").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Product Page check for hint date with invalid format"
$ws.Range("B10").Value = "<HINT Must be a valid date>"

# --- Row 11: "Product Page enter date with invalid value in past" ---
$ws.Range("B11").Value = "01/01/2000"

# --- Row 12: "Product Page check for hint date with invalid value in past" ---
$ws.Range("A6").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Product Page check for hint date with invalid value in past"
$ws.Range("B12").Value = "<HINT Must be more than one month in the future>"

# --- Row 9 / 11 labels (added after the strings above so the shared-string
#     table gets the same insertion order as the authored workbook) ---
$ws.Range("A9").Value = "Product Page enter date with invalid format"
$ws.Range("A11").Value = "Product Page enter date with invalid value in past"

# --- Move the product-screenshot picture down to sit below the new rows ---
$shp = $ws.Shapes.Item(1)
$shp.Top = 195.62188976377954
$shp.Width = 836.4

# --- Restore the selection to B1 ---
$ws.Range("B1").Select()
